# Weekly update: a new sampling date (and its market-day figures) was
# recorded for "Puerro" at Vega Modelo de Temuco. This inserts a new row
# right above the existing row 56, pushing the rest of the data block
# down by one row (old row 56 -> 57, ..., old row 131 -> 132), and fills
# the newly inserted row with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (shifts rows 56..131 down to 57..132,
# carrying over the existing column formatting, e.g. the date style on D).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new week's record.
$ws.Cells.Item(56, 1).Value  = 10
$ws.Cells.Item(56, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value  = "La Araucanía"
$ws.Cells.Item(56, 4).Value  = 44467
$ws.Cells.Item(56, 5).Value  = 9
$ws.Cells.Item(56, 6).Value  = 100112005
$ws.Cells.Item(56, 7).Value  = "Puerro"
$ws.Cells.Item(56, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(56, 9).Value  = "Primera"
$ws.Cells.Item(56, 10).Value = 70
$ws.Cells.Item(56, 11).Value = 6000
$ws.Cells.Item(56, 12).Value = 7000
$ws.Cells.Item(56, 13).Value = 6429
$ws.Cells.Item(56, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(56, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(56, 16).Value = 536
$ws.Cells.Item(56, 17).Value = 12
$ws.Cells.Item(56, 18).Value = "Hortaliza"
